# The bot's working folder moved from "C:\Projetos\bot whatsapp\" to
# "R:\#Prontos\bot_whatsapp\" on the host machine. The attachment cell
# (WhatsApp!C9) stores the absolute path to the file the bot attaches
# before it looks for the "enviar" (send) button - with the old, now
# wrong, path the dialog navigation failed and the send button could
# not be located. Point the attachment path at the new folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WhatsApp")

$ws.Range("C9").Value = "R:\#Prontos\bot_whatsapp\test.jpg"
